$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new price reading is inserted at the top of the series (row 50)
# and all the subsequent historical readings (rows 51-94) shift down by one row,
# so the former row 94 becomes the new row 95.

# 1) Capture the current (pre-edit) values for the columns that vary per reading:
#    D = Fecha, J = Volumen, K = Precio minimo, L = Precio maximo,
#    M = Precio promedio ponderado, P = Precio $/Kg
$origD = @{}
$origJ = @{}
$origK = @{}
$origL = @{}
$origM = @{}
$origP = @{}
for ($r = 50; $r -le 94; $r++) {
    $origD[$r] = $ws.Cells.Item($r, 4).Value()
    $origJ[$r] = $ws.Cells.Item($r, 10).Value()
    $origK[$r] = $ws.Cells.Item($r, 11).Value()
    $origL[$r] = $ws.Cells.Item($r, 12).Value()
    $origM[$r] = $ws.Cells.Item($r, 13).Value()
    $origP[$r] = $ws.Cells.Item($r, 16).Value()
}

# 2) Shift rows 51-94 down by one (row r takes the captured values of row r-1).
#    Going from the bottom up is safe since all source values were captured above.
for ($r = 94; $r -ge 51; $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 4).Value = $origD[$src]
    $ws.Cells.Item($r, 10).Value = $origJ[$src]
    $ws.Cells.Item($r, 11).Value = $origK[$src]
    $ws.Cells.Item($r, 12).Value = $origL[$src]
    $ws.Cells.Item($r, 13).Value = $origM[$src]
    $ws.Cells.Item($r, 16).Value = $origP[$src]
}

# 3) Row 50 holds the brand-new weekly reading (new date & volume; prices unchanged).
$ws.Cells.Item(50, 4).Value = 45216
$ws.Cells.Item(50, 10).Value = 30

# 4) Append new row 95: a full copy of the original row 94, now at the bottom of the series.
$ws.Cells.Item(95, 1).Value = 4
$ws.Cells.Item(95, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(95, 3).Value = 'Los Lagos'
$ws.Cells.Item(95, 4).Value = 44789
$ws.Cells.Item(95, 5).Value = 10
$ws.Cells.Item(95, 6).Value = 100112012
$ws.Cells.Item(95, 7).Value = 'Espinaca'
$ws.Cells.Item(95, 8).Value = 'Sin especificar'
$ws.Cells.Item(95, 9).Value = 'Primera'
$ws.Cells.Item(95, 10).Value = 30
$ws.Cells.Item(95, 11).Value = 15000
$ws.Cells.Item(95, 12).Value = 15000
$ws.Cells.Item(95, 13).Value = 15000
$ws.Cells.Item(95, 14).Value = '$/cuna 10 kilos'
$ws.Cells.Item(95, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(95, 16).Value = 1500
$ws.Cells.Item(95, 17).Value = 10
$ws.Cells.Item(95, 18).Value = 'Hortaliza'

# 5) Make sure the new row's date cell uses the same number format as the rest of column D.
$ws.Cells.Item(95, 4).NumberFormat = $ws.Cells.Item(94, 4).NumberFormat

